$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

function Copy-Format($srcRef, $dstRef) {
    $ws.Range($srcRef).Copy()
    $ws.Range($dstRef).PasteSpecial(-4122)
}

# ===== Pass 1: apply cell formatting (styles) to the new range =====
# Row 183
Copy-Format "A3" "A183"
Copy-Format "A3" "B183"
Copy-Format "A3" "C183"
Copy-Format "A3" "D183"
Copy-Format "A3" "E183"

# Row 184
Copy-Format "A21" "A184"
Copy-Format "B21" "B184"
Copy-Format "B21" "C184"
Copy-Format "B21" "D184"
Copy-Format "E21" "E184"

# Row 185
Copy-Format "A51" "A185"
Copy-Format "B9" "B185"
Copy-Format "D9" "C185"
Copy-Format "D9" "D185"
Copy-Format "A9" "E185"

# Row 186
Copy-Format "A21" "A186"
Copy-Format "B22" "B186"
Copy-Format "B21" "C186"
Copy-Format "B21" "D186"
Copy-Format "E21" "E186"

# Row 187
Copy-Format "A9" "A187"
Copy-Format "B9" "B187"
Copy-Format "D9" "C187"
Copy-Format "D9" "D187"
Copy-Format "E168" "E187"

# Row 188
Copy-Format "E21" "A188"
Copy-Format "B22" "B188"
Copy-Format "B21" "C188"
Copy-Format "B21" "D188"
Copy-Format "E21" "E188"

# Row 189
Copy-Format "A15" "A189"
Copy-Format "B15" "B189"
Copy-Format "D15" "C189"
Copy-Format "D15" "D189"
Copy-Format "A15" "E189"

# Row 190
Copy-Format "E21" "A190"
Copy-Format "B22" "B190"
Copy-Format "B21" "C190"
Copy-Format "B21" "D190"
Copy-Format "E21" "E190"

# Row 191
Copy-Format "A9" "A191"
Copy-Format "B9" "B191"
Copy-Format "D9" "C191"
Copy-Format "D9" "D191"
Copy-Format "A9" "E191"

# Row 192
Copy-Format "E21" "A192"
Copy-Format "B22" "B192"
Copy-Format "B21" "C192"
Copy-Format "B21" "D192"
Copy-Format "E21" "E192"

# Row 193
Copy-Format "A9" "A193"
Copy-Format "B9" "B193"
Copy-Format "D9" "C193"
Copy-Format "D9" "D193"
Copy-Format "A9" "E193"

# Row 194
Copy-Format "E21" "A194"
Copy-Format "B22" "B194"
Copy-Format "B21" "C194"
Copy-Format "B21" "D194"
Copy-Format "E21" "E194"

# Row 195
Copy-Format "E21" "A195"
Copy-Format "B22" "B195"
Copy-Format "B22" "C195"
Copy-Format "B21" "D195"
Copy-Format "E21" "E195"

# Row 196
Copy-Format "A9" "A196"
Copy-Format "B9" "B196"
Copy-Format "B9" "C196"
Copy-Format "D9" "D196"
Copy-Format "A9" "E196"

# Row 197
Copy-Format "A14" "A197"
Copy-Format "B14" "B197"
Copy-Format "B14" "C197"
Copy-Format "D14" "D197"
Copy-Format "A14" "E197"

# Row 198
Copy-Format "A15" "A198"
Copy-Format "B15" "B198"
Copy-Format "B15" "C198"
Copy-Format "D15" "D198"
Copy-Format "A15" "E198"

# Row 199
Copy-Format "A14" "A199"
Copy-Format "B14" "B199"
Copy-Format "B14" "C199"
Copy-Format "D14" "D199"
Copy-Format "A14" "E199"

# Row 200
Copy-Format "A29" "A200"
Copy-Format "A29" "B200"
Copy-Format "A29" "C200"
Copy-Format "D29" "D200"
Copy-Format "A29" "E200"

# Row 203
Copy-Format "A3" "A203"
Copy-Format "A3" "B203"
Copy-Format "A3" "C203"
Copy-Format "A3" "D203"
Copy-Format "A3" "E203"

# Row 204
Copy-Format "A21" "A204"
Copy-Format "B21" "B204"
Copy-Format "B21" "C204"
Copy-Format "B21" "D204"
Copy-Format "E21" "E204"

# Row 205
Copy-Format "A51" "A205"
Copy-Format "B9" "B205"
Copy-Format "D9" "C205"
Copy-Format "D9" "D205"
Copy-Format "A9" "E205"

# Row 206
Copy-Format "A21" "A206"
Copy-Format "B22" "B206"
Copy-Format "B21" "C206"
Copy-Format "B21" "D206"
Copy-Format "E21" "E206"

# Row 207
Copy-Format "E21" "A207"
Copy-Format "B22" "B207"
Copy-Format "B21" "C207"
Copy-Format "B21" "D207"
Copy-Format "E136" "E207"

# Row 208
Copy-Format "E21" "A208"
Copy-Format "B22" "B208"
Copy-Format "B21" "C208"
Copy-Format "B21" "D208"
Copy-Format "E21" "E208"

# Row 209
Copy-Format "E21" "A209"
Copy-Format "B22" "B209"
Copy-Format "B21" "C209"
Copy-Format "B21" "D209"
Copy-Format "E21" "E209"

# Row 210
Copy-Format "E21" "A210"
Copy-Format "B22" "B210"
Copy-Format "B21" "C210"
Copy-Format "B21" "D210"
Copy-Format "E21" "E210"

# Row 211
Copy-Format "E21" "A211"
Copy-Format "B22" "B211"
Copy-Format "B21" "C211"
Copy-Format "B21" "D211"
Copy-Format "E21" "E211"

# Row 212
Copy-Format "E21" "A212"
Copy-Format "B22" "B212"
Copy-Format "B21" "C212"
Copy-Format "B21" "D212"
Copy-Format "E21" "E212"

# Row 213
Copy-Format "E21" "A213"
Copy-Format "B22" "B213"
Copy-Format "B21" "C213"
Copy-Format "B21" "D213"
Copy-Format "E21" "E213"

# Row 214
Copy-Format "E21" "A214"
Copy-Format "B22" "B214"
Copy-Format "B22" "C214"
Copy-Format "B21" "D214"
Copy-Format "E21" "E214"

# Row 215
Copy-Format "A9" "A215"
Copy-Format "B9" "B215"
Copy-Format "B9" "C215"
Copy-Format "D9" "D215"
Copy-Format "A9" "E215"

# Row 216
Copy-Format "A14" "A216"
Copy-Format "B14" "B216"
Copy-Format "B14" "C216"
Copy-Format "D14" "D216"
Copy-Format "A14" "E216"

# Row 217
Copy-Format "A15" "A217"
Copy-Format "B15" "B217"
Copy-Format "B15" "C217"
Copy-Format "D15" "D217"
Copy-Format "A15" "E217"

# Row 218
Copy-Format "A14" "A218"
Copy-Format "B14" "B218"
Copy-Format "B14" "C218"
Copy-Format "D14" "D218"
Copy-Format "A14" "E218"

# Row 219
Copy-Format "A29" "A219"
Copy-Format "A29" "B219"
Copy-Format "A29" "C219"
Copy-Format "D29" "D219"
Copy-Format "A29" "E219"

# ===== Pass 2: set values / formulas =====
# Row 183
$ws.Range("A183").Value = "Datum"
$ws.Range("B183").Value = "Start"
$ws.Range("C183").Value = "Slut"
$ws.Range("D183").Value = "Längd"
$ws.Range("E183").Value = "Uppgift"

# Row 184
$ws.Range("A184").Value = 44643
$ws.Range("B184").Value = 0.38541666666666669
$ws.Range("C184").Value = 0.41666666666666669
$ws.Range("D184").Formula = "=SUM(C184-B184)"
$ws.Range("E184").Value = "Pluralsight course: Testing xUnit"

# Row 185
$ws.Range("B185").Formula = "=C184"
$ws.Range("C185").Value = 0.4826388888888889
$ws.Range("D185").Formula = "=SUM(C185-B185)"
$ws.Range("E185").Value = "Paus, Verisure m.m."

# Row 186
$ws.Range("B186").Formula = "=C185"
$ws.Range("C186").Value = 0.52083333333333337
$ws.Range("D186").Formula = "=SUM(C186-B186)"
$ws.Range("E186").Value = "Pluralsight course: Testing xUnit"

# Row 187
$ws.Range("B187").Formula = "=C186"
$ws.Range("C187").Value = 0.58333333333333337
$ws.Range("D187").Formula = "=SUM(C187-B187)"
$ws.Range("E187").Value = "Lunch"

# Row 188
$ws.Range("B188").Formula = "=C187"
$ws.Range("C188").Value = 0.63194444444444442
$ws.Range("D188").Formula = "=SUM(C188-B188)"
$ws.Range("E188").Value = "Pluralsight course: Testing xUnit"

# Row 189
$ws.Range("B189").Formula = "=C188"
$ws.Range("C189").Value = 0.65277777777777779
$ws.Range("D189").Formula = "=SUM(C189-B189)"
$ws.Range("E189").Value = "Jobbsök"

# Row 190
$ws.Range("B190").Formula = "=C189"
$ws.Range("C190").Value = 0.70833333333333304
$ws.Range("D190").Formula = "=SUM(C190-B190)"
$ws.Range("E190").Value = "Pluralsight course: Testing xUnit, UserManagerClass, UML"

# Row 191
$ws.Range("B191").Formula = "=C190"
$ws.Range("C191").Value = 0.75
$ws.Range("D191").Formula = "=SUM(C191-B191)"
$ws.Range("E191").Value = "Paus"

# Row 192
$ws.Range("B192").Formula = "=C191"
$ws.Range("C192").Value = 0.79166666666666696
$ws.Range("D192").Formula = "=SUM(C192-B192)"
$ws.Range("E192").Value = "Rnd Players"

# Row 193
$ws.Range("B193").Formula = "=C192"
$ws.Range("C193").Value = 0.86458333333333337
$ws.Range("D193").Formula = "=SUM(C193-B193)"
$ws.Range("E193").Value = "Middag"

# Row 194
$ws.Range("B194").Formula = "=C193"
$ws.Range("C194").Value = 0.89583333333333337
$ws.Range("D194").Formula = "=SUM(C194-B194)"
$ws.Range("E194").Value = "Init and List Team Players"

# Row 195

# Row 196
$ws.Range("A196").Value = "Misc Total"
$ws.Range("D196").Formula = "=D185+D187+D191+D193"

# Row 197

# Row 198
$ws.Range("A198").Value = "Jobbsök aktiviteter"
$ws.Range("D198").Formula = "=D189"

# Row 199

# Row 200
$ws.Range("A200").Value = "Total dev"
$ws.Range("D200").Formula = "=SUM(D184:D194)-(D196+D198)"

# Row 203
$ws.Range("A203").Value = "Datum"
$ws.Range("B203").Value = "Start"
$ws.Range("C203").Value = "Slut"
$ws.Range("D203").Value = "Längd"
$ws.Range("E203").Value = "Uppgift"

# Row 204
$ws.Range("A204").Value = 44644
$ws.Range("B204").Value = 0.33333333333333331
$ws.Range("C204").Value = 0.47916666666666669
$ws.Range("D204").Formula = "=SUM(C204-B204)"

# Row 205
$ws.Range("B205").Formula = "=C204"
$ws.Range("C205").Value = 0.52083333333333337
$ws.Range("D205").Formula = "=SUM(C205-B205)"
$ws.Range("E205").Value = "Paus"

# Row 206
$ws.Range("B206").Formula = "=C205"
$ws.Range("C206").Value = 0.77083333333333337
$ws.Range("D206").Formula = "=SUM(C206-B206)"

# Row 207
$ws.Range("B207").Formula = "=C206"
$ws.Range("C207").Value = 0.875
$ws.Range("D207").Formula = "=SUM(C207-B207)"

# Row 208

# Row 209

# Row 210

# Row 211

# Row 212

# Row 213

# Row 214

# Row 215
$ws.Range("A215").Value = "Misc Total"
$ws.Range("D215").Formula = "=D205"

# Row 216

# Row 217
$ws.Range("A217").Value = "Jobbsök aktiviteter"

# Row 218

# Row 219
$ws.Range("A219").Value = "Total dev"
$ws.Range("D219").Formula = "=SUM(D204:D213)-(D215+D217)"

# ===== Pass 3: row heights for thick border rows =====
$ws.Rows.Item(200).RowHeight = 15.75
$ws.Rows.Item(201).RowHeight = 15.75
$ws.Rows.Item(219).RowHeight = 15.75
$ws.Rows.Item(220).RowHeight = 15.75

# ===== Pass 4: sheet view selection =====
$ws.Range("C206").Select()

# ===== Pass 5: force full recalculation so cached <v> are fresh =====
$excel.CalculateFull()